$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---- Row 2 : update existing row (ID changes 1151 -> 1148, new data columns) ----
$ws.Range("A2").Value = 1148
$ws.Range("B2").Value = "NFLX"
Set-TextCell $ws.Range("C2") "2013"
$ws.Range("D2").Value = 1.93
$ws.Range("E2").Value = 58239896.3731
Set-TextCell $ws.Range("F2") "(+)"
Set-TextCell $ws.Range("G2") "Solid"

# ---- Row 3 : new row ----
$ws.Range("A3").Value = 1149
$ws.Range("B3").Value = "NFLX"
Set-TextCell $ws.Range("C3") "2014"
$ws.Range("D3").Value = 4.44
$ws.Range("E3").Value = 60089864.8649
Set-TextCell $ws.Range("F3") "(+)"
Set-TextCell $ws.Range("G3") "Solid"

# ---- Row 4 : new row ----
$ws.Range("A4").Value = 1150
$ws.Range("B4").Value = "NFLX"
Set-TextCell $ws.Range("C4") "2015"
$ws.Range("D4").Value = 0.29
$ws.Range("E4").Value = 422900000
Set-TextCell $ws.Range("F4") "(+)"
Set-TextCell $ws.Range("G4") "Solid"
